$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TemplateDetailRow {
    param($Row, $TemplateCode, $WebColumn, $DatabaseColumn, $QueryColumn, $DataType, $Sequence, $IsPrimary, $IsShow, $IsCreate, $IsEdit)
    $ws.Cells.Item($Row, 2).Value2 = $TemplateCode
    $ws.Cells.Item($Row, 3).Value2 = $WebColumn
    $ws.Cells.Item($Row, 4).Value2 = $DatabaseColumn
    $ws.Cells.Item($Row, 5).Value2 = $QueryColumn
    $ws.Cells.Item($Row, 6).Value2 = $DataType
    $ws.Cells.Item($Row, 7).Value2 = $Sequence
    $ws.Cells.Item($Row, 8).Value2 = $IsPrimary
    $ws.Cells.Item($Row, 9).Value2 = $IsShow
    $ws.Cells.Item($Row, 10).Value2 = $IsCreate
    $ws.Cells.Item($Row, 11).Value2 = $IsEdit

    $mFormula = '="(' + "'" + '"&B' + $Row + '&"' + "','" + '"&C' + $Row + '&"' + "','" + '"&D' + $Row + '&"' + "','" + '"&E' + $Row + '&"' + "','" + '"&F' + $Row + '&"' + "','" + '"&G' + $Row + '&"' + "','" + '"&H' + $Row + '&"' + "','" + '"&I' + $Row + '&"' + "','" + '"&J' + $Row + '&"' + "','" + '"&K' + $Row + '&"' + "');" + '"'
    $ws.Cells.Item($Row, 13).Formula = $mFormula
    $ws.Cells.Item($Row, 14).Formula = '="INSERT INTO INDO_CMS_TEMPLATE_DETAIL (template_code,web_column,database_column,query_column,data_type,sequence,is_primary,is_show,is_create,is_edit) VALUES"'
    $ws.Cells.Item($Row, 16).Formula = "=N$Row" + "&M$Row"
}

# --- Step 1: Insert a new row at 60 -> INDO_CMS_TEMPLATE_DETAIL / MAX LENGTH (sequence 13) ---
$ws.Rows.Item(60).Insert()
Set-TemplateDetailRow 60 "INDO_CMS_TEMPLATE_DETAIL" "MAX LENGTH" "max_length" "max_length" "INTEGER" 13 0 1 1 1

# --- Step 2: Insert a new row at 64 -> INDO_CMS_USER_APPROVER / ROLE ID (sequence 3) ---
# (this pushes the old APPROVER USERNAME/ROLE/SEQUENCE rows down and their sequence numbers get updated below)
$ws.Rows.Item(64).Insert()
Set-TemplateDetailRow 64 "INDO_CMS_USER_APPROVER" "ROLE ID" "role_id" "role_id" "STRING" 3 0 1 1 1

# Update sequence numbers on the rows that were pushed down (APPROVER USERNAME/ROLE/SEQUENCE): 3->4, 4->5, 5->6
$ws.Cells.Item(65, 7).Value2 = 4
$ws.Cells.Item(66, 7).Value2 = 5
$ws.Cells.Item(67, 7).Value2 = 6

# --- Step 3: Add INDO_CMS_JOB_HEADER rows (69-71), leaving row 68 blank as a separator ---
Set-TemplateDetailRow 69 "INDO_CMS_JOB_HEADER" "ROW ID" "row_id" "row_id" "INTEGER" 1 1 1 0 0
Set-TemplateDetailRow 70 "INDO_CMS_JOB_HEADER" "JOB ID" "job_id" "job_id" "STRING" 2 0 1 1 1
Set-TemplateDetailRow 71 "INDO_CMS_JOB_HEADER" "JOB DESCRIPTION" "job_description" "job_description" "STRING" 3 0 1 1 1

# --- Step 4: Add INDO_CMS_JOB_DETAIL rows (73-74), leaving row 72 blank as a separator ---
Set-TemplateDetailRow 73 "INDO_CMS_JOB_DETAIL" "ROW ID" "row_id" "row_id" "INTEGER" 1 1 1 0 0
Set-TemplateDetailRow 74 "INDO_CMS_JOB_DETAIL" "JOB ID" "job_id" "job_id" "STRING" 2 0 1 1 1

# --- Adjust column E width (widened) ---
$ws.Columns.Item(5).ColumnWidth = 24.07

# --- Update sheet1 view / selection to match final state ---
$ws.Application.ActiveWindow.ScrollRow = 48
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C74").Select()

Write-Output "Sheet1 edits complete"
